$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-45: simple Price (D) / Volume(1h) (E) updates ---
$ws.Range("D2").Value = "29.616.49"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.853.00"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6314"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07489"
$ws.Range("E8").Value = "  -1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2918"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.15"
$ws.Range("E10").Value = "  +2.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07751"
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "1.851.60"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6826"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001038"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.89"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.335"
$ws.Range("E17").Value = "  +3.71%  "
$ws.Range("D18").Value = "29.603.95"
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.61"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.570"
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "159.70"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.525"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1369"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.60"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06595"
$ws.Range("E28").Value = "  +16.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.462"
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.487"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.116"
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.115"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.848"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6997"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.569"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01868"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("D38").Value = "1.263.39"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.845"
$ws.Range("E39").Value = "  +4.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.785"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9384"
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("D42").Value = "2.036.56"
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.43"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.32"
$ws.Range("E45").Value = "  +1.23%  "

# --- Rows 46-51: a new coin (BabyDogeCoin) was inserted at row 46, ---
# --- shifting RenderToken..TheSandbox down one row each; Cronos drops off ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000120"
$ws.Range("E46").Value = "  +4.72%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.744"
$ws.Range("E47").Value = "  +4.10%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.128"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1162"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.027"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3957"
$ws.Range("E51").Value = "  -0.98%  "
